# Weekly update: insert a new "Cebollín baby" price record at row 52,
# pushing the existing rows 52-104 down to 53-105 (dimension grows to A1:R105).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 52:104 down by one row, opening up a blank row 52.
$ws.Rows("52:52").Insert()

# Populate the newly inserted row 52 with the new weekly record.
$ws.Range("A52").Value = 1
$ws.Range("B52").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C52").Value = "Arica y Parinacota"
$ws.Range("D52").Value = 44827
$ws.Range("E52").Value = 15
$ws.Range("F52").Value = 100112038
$ws.Range("G52").Value = "Cebollín baby"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 300
$ws.Range("K52").Value = 900
$ws.Range("L52").Value = 1000
$ws.Range("M52").Value = 950
$ws.Range("N52").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O52").Value = "Región de Arica y Parinacota"
$ws.Range("P52").Value = 475
$ws.Range("Q52").Value = 2
$ws.Range("R52").Value = "Hortaliza"

# Match the date format used by the other rows in column D.
$ws.Range("D52").NumberFormat = $ws.Range("D53").NumberFormat
